$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arr = New-Object "object[,]" 24,12
$arr[0,0] = 14.73660810783698
$arr[0,1] = 17.80111857467285
$arr[0,2] = 14.6057720875824
$arr[0,3] = 15.27805400130843
$arr[0,4] = 0
$arr[0,5] = 3.831539936327273
$arr[0,6] = 0
$arr[0,7] = 54.94166074240064
$arr[0,8] = 8.9959315605343
$arr[0,9] = 18.05490431119427
$arr[0,10] = 0
$arr[0,11] = 21.58195125715804
$arr[1,0] = 14.8225020206052
$arr[1,1] = 17.57506533151998
$arr[1,2] = 14.58691156789853
$arr[1,3] = 15.27395774500216
$arr[1,4] = 0
$arr[1,5] = 3.835835349899482
$arr[1,6] = 0
$arr[1,7] = 53.64733881347002
$arr[1,8] = 9.010161120887702
$arr[1,9] = 18.07295164217562
$arr[1,10] = 0
$arr[1,11] = 21.61392487222412
$arr[2,0] = 14.88245211535165
$arr[2,1] = 17.43796371365347
$arr[2,2] = 14.57821892355886
$arr[2,3] = 15.27399034997376
$arr[2,4] = 0
$arr[2,5] = 3.838602358811631
$arr[2,6] = 0
$arr[2,7] = 52.83364144276865
$arr[2,8] = 9.019642811693082
$arr[2,9] = 18.08997821030006
$arr[2,10] = 0
$arr[2,11] = 21.63866654137562
$arr[3,0] = 14.9086800710223
$arr[3,1] = 17.3825731461305
$arr[3,2] = 14.5754040382008
$arr[3,3] = 15.27464444660636
$arr[3,4] = 0
$arr[3,5] = 3.839762674575924
$arr[3,6] = 0
$arr[3,7] = 52.49757089651473
$arr[3,8] = 9.02369424822365
$arr[3,9] = 18.09840654971507
$arr[3,10] = 0
$arr[3,11] = 21.65003082234906
$arr[4,0] = 14.91314332889524
$arr[4,1] = 17.37340607132736
$arr[4,2] = 14.57498059201098
$arr[4,3] = 15.27479176195194
$arr[4,4] = 0
$arr[4,5] = 3.839957325652892
$arr[4,6] = 0
$arr[4,7] = 52.44150470861128
$arr[4,8] = 9.024378324204335
$arr[4,9] = 18.09989588627494
$arr[4,10] = 0
$arr[4,11] = 21.65199517705834
$arr[5,0] = 14.88279857604717
$arr[5,1] = 17.43721468650936
$arr[5,2] = 14.57817801430897
$arr[5,3] = 15.27399657668932
$arr[5,4] = 0
$arr[5,5] = 3.838617874489446
$arr[5,6] = 0
$arr[5,7] = 52.82912683152864
$arr[5,8] = 9.019696690893223
$arr[5,9] = 18.09008585278675
$arr[5,10] = 0
$arr[5,11] = 21.63881461812645
$arr[6,0] = 14.76471852762851
$arr[6,1] = 17.72285547005847
$arr[6,2] = 14.59867052413581
$arr[6,3] = 15.27611300055153
$arr[6,4] = 0
$arr[6,5] = 3.832994183228483
$arr[6,6] = 0
$arr[6,7] = 54.49948645804034
$arr[6,8] = 9.0006835421855
$arr[6,9] = 18.0598902132853
$arr[6,10] = 0
$arr[6,11] = 21.59191351308058
$arr[7,0] = 14.59105286522638
$arr[7,1] = 18.29409170414871
$arr[7,2] = 14.6616950820466
$arr[7,3] = 15.30046319357892
$arr[7,4] = 0
$arr[7,5] = 3.822987825435478
$arr[7,6] = 0
$arr[7,7] = 57.6125913878368
$arr[7,8] = 8.969293785712933
$arr[7,9] = 18.04804323295177
$arr[7,10] = 0
$arr[7,11] = 21.54060627327451
$arr[8,0] = 14.49964564075593
$arr[8,1] = 18.7173622321192
$arr[8,2] = 14.72179598491216
$arr[8,3] = 15.33063572259131
$arr[8,4] = 0
$arr[8,5] = 3.81624959684214
$arr[8,6] = 0
$arr[8,7] = 59.78609433121997
$arr[8,8] = 8.949807021203767
$arr[8,9] = 18.06842127093461
$arr[8,10] = 0
$arr[8,11] = 21.52784862882674
$arr[9,0] = 14.46610071679534
$arr[9,1] = 18.9100251736756
$arr[9,2] = 14.75209641324212
$arr[9,3] = 15.34701685218419
$arr[9,4] = 0
$arr[9,5] = 3.813315382978558
$arr[9,6] = 0
$arr[9,7] = 60.74744612116969
$arr[9,8] = 8.941714461871761
$arr[9,9] = 18.08403080280786
$arr[9,10] = 0
$arr[9,11] = 21.52748293553594
$arr[10,0] = 14.45456780640485
$arr[10,1] = 18.98294277679249
$arr[10,2] = 14.76399208380897
$arr[10,3] = 15.35360023993409
$arr[10,4] = 0
$arr[10,5] = 3.81222295883554
$arr[10,6] = 0
$arr[10,7] = 61.10734104043758
$arr[10,8] = 8.938760737466387
$arr[10,9] = 18.09085360131148
$arr[10,10] = 0
$arr[10,11] = 21.52812748891343
$arr[11,0] = 14.45699937851531
$arr[11,1] = 18.96724131312248
$arr[11,2] = 14.76141146017135
$arr[11,3] = 15.35216550886141
$arr[11,4] = 0
$arr[11,5] = 3.812457402579772
$arr[11,6] = 0
$arr[11,7] = 61.03001887854072
$arr[11,8] = 8.939391953414995
$arr[11,9] = 18.0893436444114
$arr[11,10] = 0
$arr[11,11] = 21.52795383483643
$arr[12,0] = 14.46512837362338
$arr[12,1] = 18.91602526162562
$arr[12,2] = 14.75306665390166
$arr[12,3] = 15.3475508619143
$arr[12,4] = 0
$arr[12,5] = 3.813225134660191
$arr[12,6] = 0
$arr[12,7] = 60.77713916866955
$arr[12,8] = 8.9414692390902
$arr[12,9] = 18.08457384707949
$arr[12,10] = 0
$arr[12,11] = 21.52752026561631
$arr[13,0] = 14.47026037213482
$arr[13,1] = 18.8846471364578
$arr[13,2] = 14.74800998771793
$arr[13,3] = 15.34477371828029
$arr[13,4] = 0
$arr[13,5] = 3.813697823594467
$arr[13,6] = 0
$arr[13,7] = 60.62169689894134
$arr[13,8] = 8.942756051760586
$arr[13,9] = 18.0817709411311
$arr[13,10] = 0
$arr[13,11] = 21.52735669027013
$arr[14,0] = 14.50200082462366
$arr[14,1] = 18.70476873964934
$arr[14,2] = 14.71987500662601
$arr[14,3] = 15.32961847812203
$arr[14,4] = 0
$arr[14,5] = 3.816443979135713
$arr[14,6] = 0
$arr[14,7] = 59.7226982973247
$arr[14,8] = 8.950351401823198
$arr[14,9] = 18.06752875811181
$arr[14,10] = 0
$arr[14,11] = 21.52798205351605
$arr[15,0] = 14.52354123483004
$arr[15,1] = 18.59440829208323
$arr[15,2] = 14.7033705121685
$arr[15,3] = 15.32100028139403
$arr[15,4] = 0
$arr[15,5] = 3.818162116389079
$arr[15,6] = 0
$arr[15,7] = 59.16402254227088
$arr[15,8] = 8.95520845526525
$arr[15,9] = 18.06041579704414
$arr[15,10] = 0
$arr[15,11] = 21.52975936478266
$arr[16,0] = 14.53668635290445
$arr[16,1] = 18.53094601319058
$arr[16,2] = 14.69415647146072
$arr[16,3] = 15.3162933867226
$arr[16,4] = 0
$arr[16,5] = 3.819162686256836
$arr[16,6] = 0
$arr[16,7] = 58.84012079020401
$arr[16,8] = 8.958074790658403
$arr[16,9] = 18.05692136333429
$arr[16,10] = 0
$arr[16,11] = 21.53129342522931
$arr[17,0] = 14.54126640413729
$arr[17,1] = 18.5094629342639
$arr[17,2] = 14.69108477601201
$arr[17,3] = 15.31474270743852
$arr[17,4] = 0
$arr[17,5] = 3.819503586299638
$arr[17,6] = 0
$arr[17,7] = 58.73001916144658
$arr[17,8] = 8.959057774839447
$arr[17,9] = 18.0558406780322
$arr[17,10] = 0
$arr[17,11] = 21.53190069303222
$arr[18,0] = 14.52116991045137
$arr[18,1] = 18.6061552984098
$arr[18,2] = 14.70509860620226
$arr[18,3] = 15.32189183351422
$arr[18,4] = 0
$arr[18,5] = 3.817977941363312
$arr[18,6] = 0
$arr[18,7] = 59.2237617036967
$arr[18,8] = 8.954683892800839
$arr[18,9] = 18.06111121748063
$arr[18,10] = 0
$arr[18,11] = 21.5295171885424
$arr[19,0] = 14.46270883253489
$arr[19,1] = 18.93107017626433
$arr[19,2] = 14.75550632007549
$arr[19,3] = 15.34889599131391
$arr[19,4] = 0
$arr[19,5] = 3.812999126804713
$arr[19,6] = 0
$arr[19,7] = 60.85153023598049
$arr[19,8] = 8.940856086294996
$arr[19,9] = 18.08595010947607
$arr[19,10] = 0
$arr[19,11] = 21.52762635775692
$arr[20,0] = 14.43132584370215
$arr[20,1] = 19.14316441519523
$arr[20,2] = 14.79090558195781
$arr[20,3] = 15.36875996186812
$arr[20,4] = 0
$arr[20,5] = 3.809854108417243
$arr[20,6] = 0
$arr[20,7] = 61.89110739304969
$arr[20,8] = 8.932464239757946
$arr[20,9] = 18.10749783971417
$arr[20,10] = 0
$arr[20,11] = 21.53095488522191
$arr[21,0] = 14.44744657544791
$arr[21,1] = 19.03000744300189
$arr[21,2] = 14.77178918900435
$arr[21,3] = 15.35795611250799
$arr[21,4] = 0
$arr[21,5] = 3.811522745058037
$arr[21,6] = 0
$arr[21,7] = 61.33855051912083
$arr[21,8] = 8.936884158605766
$arr[21,9] = 18.09551137150464
$arr[21,10] = 0
$arr[21,11] = 21.52876051577627
$arr[22,0] = 14.52223961616249
$arr[22,1] = 18.60084451800927
$arr[22,2] = 14.7043164792246
$arr[22,3] = 15.32148799077058
$arr[22,4] = 0
$arr[22,5] = 3.818061167043866
$arr[22,6] = 0
$arr[22,7] = 59.19676205876989
$arr[22,8] = 8.954920817130557
$arr[22,9] = 18.06079496494055
$arr[22,10] = 0
$arr[22,11] = 21.52962508093056
$arr[23,0] = 14.63174248296708
$arr[23,1] = 18.13869779219885
$arr[23,2] = 14.64221130181662
$arr[23,3] = 15.29171549646489
$arr[23,4] = 0
$arr[23,5] = 3.825586403996957
$arr[23,6] = 0
$arr[23,7] = 56.78938856588738
$arr[23,8] = 8.977156362341129
$arr[23,9] = 18.0461521686599
$arr[23,10] = 0
$arr[23,11] = 21.55011588052989

$ws.Range("B2:M25").Value = $arr
